$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 45.1
$ws.Cells.Item(2,7).Value = 110.2
$ws.Cells.Item(2,8).Value = 55.5
$ws.Cells.Item(2,9).Value = 40.7
$ws.Cells.Item(3,3).Value = 45.1
$ws.Cells.Item(3,7).Value = 245.2
$ws.Cells.Item(3,8).Value = 55.5
$ws.Cells.Item(4,3).Value = 45.1
$ws.Cells.Item(4,7).Value = 35.3
$ws.Cells.Item(4,8).Value = 55.5
$ws.Cells.Item(4,9).Value = 33.5
$ws.Cells.Item(5,3).Value = 45.1
$ws.Cells.Item(5,7).Value = 105.2
$ws.Cells.Item(5,8).Value = 55.5
$ws.Cells.Item(5,9).Value = 40.7
$ws.Cells.Item(6,3).Value = 45.1
$ws.Cells.Item(6,7).Value = 170.07
$ws.Cells.Item(6,8).Value = 55.5
$ws.Cells.Item(7,3).Value = 45.1
$ws.Cells.Item(7,7).Value = 8.09
$ws.Cells.Item(7,8).Value = 55.5
$ws.Cells.Item(7,9).Value = 33.5
$ws.Cells.Item(8,3).Value = 45.1
$ws.Cells.Item(8,7).Value = 162.12
$ws.Cells.Item(8,8).Value = 55.5
$ws.Cells.Item(8,9).Value = 40.7
$ws.Cells.Item(9,3).Value = 45.1
$ws.Cells.Item(9,7).Value = 141.23
$ws.Cells.Item(9,8).Value = 55.5
$ws.Cells.Item(10,3).Value = 45.1
$ws.Cells.Item(10,7).Value = 22.7
$ws.Cells.Item(10,8).Value = 55.5
$ws.Cells.Item(10,9).Value = 33.5
$ws.Cells.Item(11,3).Value = 45.1
$ws.Cells.Item(11,7).Value = 16.34
$ws.Cells.Item(11,8).Value = 55.5
$ws.Cells.Item(11,9).Value = 40.7
$ws.Cells.Item(12,3).Value = 45.1
$ws.Cells.Item(12,7).Value = 95.09
$ws.Cells.Item(12,8).Value = 55.5
$ws.Cells.Item(13,3).Value = 45.1
$ws.Cells.Item(13,7).Value = 64.4
$ws.Cells.Item(13,8).Value = 55.5
$ws.Cells.Item(13,9).Value = 33.5
$ws.Cells.Item(14,3).Value = 45.1
$ws.Cells.Item(14,7).Value = 32.67
$ws.Cells.Item(14,8).Value = 55.5
$ws.Cells.Item(14,9).Value = 40.7
$ws.Cells.Item(15,3).Value = 45.1
$ws.Cells.Item(15,7).Value = 12.41
$ws.Cells.Item(15,8).Value = 55.5
$ws.Cells.Item(16,3).Value = 45.1
$ws.Cells.Item(16,7).Value = 23.12
$ws.Cells.Item(16,8).Value = 55.5
$ws.Cells.Item(16,9).Value = 33.5
$ws.Cells.Item(17,3).Value = 45.1
$ws.Cells.Item(17,7).Value = 15.07
$ws.Cells.Item(17,8).Value = 55.5
$ws.Cells.Item(17,9).Value = 40.7
$ws.Cells.Item(18,3).Value = 45.1
$ws.Cells.Item(18,7).Value = 70.2
$ws.Cells.Item(18,8).Value = 55.5
$ws.Cells.Item(19,3).Value = 45.1
$ws.Cells.Item(19,7).Value = 55.2
$ws.Cells.Item(19,8).Value = 55.5
$ws.Cells.Item(19,9).Value = 33.5
$ws.Cells.Item(20,3).Value = 45.1
$ws.Cells.Item(20,7).Value = 32.89
$ws.Cells.Item(20,8).Value = 55.5
$ws.Cells.Item(20,9).Value = 40.7
$ws.Cells.Item(21,3).Value = 45.1
$ws.Cells.Item(21,6).Value = 4
$ws.Cells.Item(21,7).Value = 20.2
$ws.Cells.Item(21,8).Value = 55.5
$ws.Cells.Item(22,3).Value = 45.1
$ws.Cells.Item(22,6).Value = 5
$ws.Cells.Item(22,7).Value = 45.2
$ws.Cells.Item(22,8).Value = 55.5
$ws.Cells.Item(22,9).Value = 33.5
$ws.Cells.Item(23,3).Value = 45.1
$ws.Cells.Item(23,6).Value = 5
$ws.Cells.Item(23,7).Value = 35.3
$ws.Cells.Item(23,8).Value = 55.5
$ws.Cells.Item(23,9).Value = 40.7
$ws.Cells.Item(24,3).Value = 45.1
$ws.Cells.Item(24,6).Value = 12
$ws.Cells.Item(24,7).Value = 15.2
$ws.Cells.Item(24,8).Value = 55.5
$ws.Cells.Item(25,3).Value = 45.1
$ws.Cells.Item(25,6).Value = 12
$ws.Cells.Item(25,7).Value = 17.07
$ws.Cells.Item(25,8).Value = 55.5
$ws.Cells.Item(25,9).Value = 33.5
$ws.Cells.Item(26,3).Value = 45.1
$ws.Cells.Item(26,6).Value = 9
$ws.Cells.Item(26,7).Value = 80.9
$ws.Cells.Item(26,8).Value = 55.5
$ws.Cells.Item(26,9).Value = 40.7
$ws.Cells.Item(27,3).Value = 45.1
$ws.Cells.Item(27,6).Value = 4
$ws.Cells.Item(27,7).Value = 16.12
$ws.Cells.Item(27,8).Value = 55.5
$ws.Cells.Item(28,3).Value = 45.1
$ws.Cells.Item(28,6).Value = 3
$ws.Cells.Item(28,7).Value = 14.23
$ws.Cells.Item(28,8).Value = 55.5
$ws.Cells.Item(28,9).Value = 33.5
$ws.Cells.Item(29,3).Value = 45.1
$ws.Cells.Item(29,6).Value = 16
$ws.Cells.Item(29,7).Value = 22.7
$ws.Cells.Item(29,8).Value = 55.5
$ws.Cells.Item(29,9).Value = 40.7
$ws.Cells.Item(30,3).Value = 45.1
$ws.Cells.Item(30,6).Value = 6
$ws.Cells.Item(30,7).Value = 16.34
$ws.Cells.Item(30,8).Value = 55.5
$ws.Cells.Item(31,3).Value = 45.1
$ws.Cells.Item(31,6).Value = 13
$ws.Cells.Item(31,7).Value = 45.09
$ws.Cells.Item(31,8).Value = 55.5
$ws.Cells.Item(31,9).Value = 33.5
$ws.Cells.Item(32,3).Value = 45.1
$ws.Cells.Item(32,6).Value = 8
$ws.Cells.Item(32,7).Value = 64.4
$ws.Cells.Item(32,8).Value = 55.5
$ws.Cells.Item(32,9).Value = 40.7
$ws.Cells.Item(33,3).Value = 45.1
$ws.Cells.Item(33,6).Value = 9
$ws.Cells.Item(33,7).Value = 32.67
$ws.Cells.Item(33,8).Value = 55.5
$ws.Cells.Item(34,3).Value = 45.1
$ws.Cells.Item(34,6).Value = 7
$ws.Cells.Item(34,7).Value = 124.1
$ws.Cells.Item(34,8).Value = 55.5
$ws.Cells.Item(34,9).Value = 33.5
$ws.Cells.Item(35,3).Value = 45.1
$ws.Cells.Item(35,6).Value = 8
$ws.Cells.Item(35,7).Value = 231.2
$ws.Cells.Item(35,8).Value = 55.5
$ws.Cells.Item(35,9).Value = 40.7
$ws.Cells.Item(36,3).Value = 45.1
$ws.Cells.Item(36,6).Value = 4
$ws.Cells.Item(36,7).Value = 150.7
$ws.Cells.Item(36,8).Value = 55.5
$ws.Cells.Item(37,3).Value = 45.1
$ws.Cells.Item(37,6).Value = 5
$ws.Cells.Item(37,7).Value = 73.3
$ws.Cells.Item(37,8).Value = 55.5
$ws.Cells.Item(37,9).Value = 33.5
$ws.Cells.Item(38,3).Value = 45.1
$ws.Cells.Item(38,6).Value = 2
$ws.Cells.Item(38,7).Value = 55.2
$ws.Cells.Item(38,8).Value = 55.5
$ws.Cells.Item(38,9).Value = 40.7
$ws.Cells.Item(39,3).Value = 45.1
$ws.Cells.Item(39,6).Value = 3
$ws.Cells.Item(39,7).Value = 32.89
$ws.Cells.Item(39,8).Value = 55.5
$ws.Cells.Item(40,3).Value = 45.1
$ws.Cells.Item(40,6).Value = 4
$ws.Cells.Item(40,7).Value = 201.2
$ws.Cells.Item(40,8).Value = 55.5
$ws.Cells.Item(40,9).Value = 33.5
$ws.Cells.Item(41,3).Value = 45.1
$ws.Cells.Item(41,6).Value = 1
$ws.Cells.Item(41,7).Value = 145.2
$ws.Cells.Item(41,8).Value = 55.5
$ws.Cells.Item(41,9).Value = 40.7
$ws.Cells.Item(42,3).Value = 45.1
$ws.Cells.Item(42,6).Value = 3
$ws.Cells.Item(42,7).Value = 135.3
$ws.Cells.Item(42,8).Value = 55.5
$ws.Cells.Item(43,3).Value = 45.1
$ws.Cells.Item(43,6).Value = 1
$ws.Cells.Item(43,7).Value = 150.2
$ws.Cells.Item(43,8).Value = 55.5
$ws.Cells.Item(43,9).Value = 33.5
$ws.Cells.Item(44,3).Value = 45.1
$ws.Cells.Item(44,6).Value = 5
$ws.Cells.Item(44,7).Value = 17.07
$ws.Cells.Item(44,8).Value = 55.5
$ws.Cells.Item(44,9).Value = 40.7
$ws.Cells.Item(45,3).Value = 45.1
$ws.Cells.Item(45,6).Value = 8
$ws.Cells.Item(45,7).Value = 80.9
$ws.Cells.Item(45,8).Value = 55.5
$ws.Cells.Item(46,3).Value = 45.1
$ws.Cells.Item(46,6).Value = 10
$ws.Cells.Item(46,7).Value = 163.12
$ws.Cells.Item(46,8).Value = 55.5
$ws.Cells.Item(46,9).Value = 33.5
$ws.Cells.Item(47,3).Value = 45.1
$ws.Cells.Item(47,6).Value = 7
$ws.Cells.Item(47,7).Value = 104.23
$ws.Cells.Item(47,8).Value = 55.5
$ws.Cells.Item(47,9).Value = 40.7
$ws.Cells.Item(48,3).Value = 45.1
$ws.Cells.Item(48,6).Value = 11
$ws.Cells.Item(48,7).Value = 221.7
$ws.Cells.Item(48,8).Value = 55.5
$ws.Cells.Item(48,9).Value = 33.5
$ws.Cells.Item(49,3).Value = 45.1
$ws.Cells.Item(49,6).Value = 3
$ws.Cells.Item(49,7).Value = 106.34
$ws.Cells.Item(49,8).Value = 55.5
$ws.Cells.Item(49,9).Value = 33.5
$ws.Cells.Item(50,3).Value = 45.1
$ws.Cells.Item(50,6).Value = 14
$ws.Cells.Item(50,7).Value = 45.09
$ws.Cells.Item(50,8).Value = 55.5
$ws.Cells.Item(50,9).Value = 40.7
$ws.Cells.Item(51,3).Value = 45.1
$ws.Cells.Item(51,6).Value = 6
$ws.Cells.Item(51,7).Value = 64.4
$ws.Cells.Item(51,8).Value = 55.5
$ws.Cells.Item(52,3).Value = 45.1
$ws.Cells.Item(52,6).Value = 7
$ws.Cells.Item(52,7).Value = 32.67
$ws.Cells.Item(52,8).Value = 55.5
$ws.Cells.Item(52,9).Value = 33.5
$ws.Cells.Item(53,3).Value = 45.1
$ws.Cells.Item(53,6).Value = 7
$ws.Cells.Item(53,7).Value = 122.41
$ws.Cells.Item(53,8).Value = 55.5
$ws.Cells.Item(53,9).Value = 40.7
$ws.Cells.Item(54,3).Value = 45.1
$ws.Cells.Item(54,6).Value = 8
$ws.Cells.Item(54,7).Value = 203.12
$ws.Cells.Item(54,8).Value = 55.5
$ws.Cells.Item(55,3).Value = 45.1
$ws.Cells.Item(55,6).Value = 4
$ws.Cells.Item(55,7).Value = 15.07
$ws.Cells.Item(55,8).Value = 55.5
$ws.Cells.Item(55,9).Value = 33.5
$ws.Cells.Item(56,3).Value = 45.1
$ws.Cells.Item(56,6).Value = 10
$ws.Cells.Item(56,7).Value = 107.2
$ws.Cells.Item(56,8).Value = 55.5
$ws.Cells.Item(56,9).Value = 40.7
$ws.Cells.Item(57,3).Value = 45.1
$ws.Cells.Item(57,6).Value = 2
$ws.Cells.Item(57,7).Value = 155.2
$ws.Cells.Item(57,8).Value = 55.5
$ws.Cells.Item(58,3).Value = 45.1
$ws.Cells.Item(58,6).Value = 6
$ws.Cells.Item(58,7).Value = 132.89
$ws.Cells.Item(58,8).Value = 55.5
$ws.Cells.Item(58,9).Value = 33.5

$ws.Range("G3").Select()
